$wb = $excel.ActiveWorkbook

# Sheet "VENTAS POR GRUPO": D30 label "0 de 28" -> "1 de 28"
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("D30").Value = "1 de 28"

# Sheet "VENTA MENSUAL": D8 (abril) and F8 (junio) go from 0 to 91.58
# F30 (junio total) is recalculated accordingly
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("D8").Value = 91.58
$wsVentaMensual.Range("F8").Value = 91.58
$wsVentaMensual.Range("F30").Value = 3893.56

# Sheet "CUMPLIMIENTO MENSUAL": widen column F and update VENTA/POR CUMPLIR/CUMPLIMIENTO
$wsCumplimientoMensual = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# Stored XML width = ColumnWidth + 5/6, so use 24.1666... to land on an XML width of 25
$wsCumplimientoMensual.Columns.Item(6).ColumnWidth = 24.166666666666668

$wsCumplimientoMensual.Range("D3").Value = 91.58
$wsCumplimientoMensual.Range("E3").Value = 3028.5345
$wsCumplimientoMensual.Range("F3").Value = 0.02935148694062349

$wsCumplimientoMensual.Range("D19").Value = 3887.8
$wsCumplimientoMensual.Range("E19").Value = 25649.99107555787
$wsCumplimientoMensual.Range("F19").Value = 0.1316212167001581
